$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4770

$ws.Range("H28").Value = 358
$ws.Range("I28").Value = 358
$ws.Range("K28").Value = 358
$ws.Range("M28").Value = 127

$ws.Range("H33").Value = 205.2
$ws.Range("I33").Value = 164.57143
$ws.Range("K33").Value = 164.57143
$ws.Range("M33").Value = 64.42857000000001

$ws.Range("H38").Value = 627.4
$ws.Range("J38").Value = 1725.6666
$ws.Range("L38").Value = 5176.9998
$ws.Range("N38").Value = -5920.9998

$ws.Range("H39").Value = 36.214287
$ws.Range("I39").Value = 16.8
$ws.Range("J39").Value = 84.75
$ws.Range("K39").Value = 50.40000000000001
$ws.Range("L39").Value = 254.25
$ws.Range("M39").Value = 245.6
$ws.Range("N39").Value = -846.25

$ws.Range("H40").Value = 2307.6924
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -3350

$ws.Range("H51").Value = 56854.2
$ws.Range("I51").Value = 10169.857
$ws.Range("J51").Value = 81991.92
$ws.Range("K51").Value = 10169.857
$ws.Range("L51").Value = 81991.92
$ws.Range("M51").Value = -9685.857
$ws.Range("N51").Value = -82959.92

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 500744.5
$ws.Range("J86").Value = 999999
$ws.Range("L86").Value = 999999
$ws.Range("N86").Value = -1002245

$ws.Range("H89").Value = 500744.5
$ws.Range("J89").Value = 999999
$ws.Range("L89").Value = 4999995
$ws.Range("N89").Value = -5011227

$ws.Range("H131").Value = 750
$ws.Range("I131").Value = 750
$ws.Range("K131").Value = 2250
$ws.Range("M131").Value = 2790

$ws.Range("H141").Value = 3439
$ws.Range("I141").Value = 2388
$ws.Range("K141").Value = 7164
$ws.Range("M141").Value = -1984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4270.6665
$ws.Range("I32").Value = 2742.3438
$ws.Range("K32").Value = 2742.3438
$ws.Range("M32").Value = -2455.3438

$ws.Range("H44").Value = 23329.334
$ws.Range("J44").Value = 23329.334
$ws.Range("L44").Value = 23329.334
$ws.Range("N44").Value = -24305.334

$ws.Range("H61").Value = 9412.625
$ws.Range("I61").Value = 9412.625
$ws.Range("K61").Value = 9412.625
$ws.Range("M61").Value = -9200.625

$ws.Range("H110").Value = 2352.2083
$ws.Range("I110").Value = 1311.7693
$ws.Range("K110").Value = 1311.7693
$ws.Range("M110").Value = 733.2307000000001

$ws.Range("H122").Value = 18659.924
$ws.Range("I122").Value = 15919.956
$ws.Range("K122").Value = 47759.868
$ws.Range("M122").Value = -45309.868

$ws.Range("H132").Value = 3706.8462
$ws.Range("I132").Value = 3962.25
$ws.Range("J132").Value = 3298.2
$ws.Range("K132").Value = 11886.75
$ws.Range("L132").Value = 9894.599999999999
$ws.Range("M132").Value = -9356.75
$ws.Range("N132").Value = -14954.6

$ws.Range("H136").Value = 9412.625
$ws.Range("I136").Value = 9412.625
$ws.Range("K136").Value = 28237.875
$ws.Range("M136").Value = -25687.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2966.3845
$ws.Range("J16").Value = 4266
$ws.Range("L16").Value = 4266
$ws.Range("N16").Value = -4840

$ws.Range("H106").Value = 11390.333
$ws.Range("J106").Value = 11390.333
$ws.Range("L106").Value = 11390.333
$ws.Range("N106").Value = -13914.333

$ws.Range("H113").Value = 2966.3845
$ws.Range("J113").Value = 4266
$ws.Range("L113").Value = 4266
$ws.Range("N113").Value = -8606

$ws.Range("H132").Value = 1599.8334
$ws.Range("I132").Value = 1809.8
$ws.Range("J132").Value = 550
$ws.Range("K132").Value = 5429.4
$ws.Range("L132").Value = 1650
$ws.Range("M132").Value = -2899.4
$ws.Range("N132").Value = -6710

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 111.77778
$ws.Range("I6").Value = 87.28570999999999
$ws.Range("K6").Value = 261.85713
$ws.Range("M6").Value = -148.85713

$ws.Range("H68").Value = 613
$ws.Range("J68").Value = 335.5
$ws.Range("L68").Value = 1006.5
$ws.Range("N68").Value = -2628.5

$ws.Range("H71").Value = 613
$ws.Range("J71").Value = 335.5
$ws.Range("L71").Value = 3019.5
$ws.Range("N71").Value = -11131.5

$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 12000
$ws.Range("M80").Value = -11064

$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 36000
$ws.Range("M83").Value = -31320

$ws.Range("H131").Value = 1381.1063
$ws.Range("I131").Value = 1134.25
$ws.Range("J131").Value = 1404.0698
$ws.Range("K131").Value = 3402.75
$ws.Range("L131").Value = 4212.2094
$ws.Range("M131").Value = 1637.25
$ws.Range("N131").Value = -14292.2094

$ws.Range("H140").Value = 4700
$ws.Range("I140").Value = 1700.2
$ws.Range("K140").Value = 5100.6
$ws.Range("M140").Value = 79.39999999999964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4253
$ws.Range("I70").Value = 4253
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4253
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3983
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 4253
$ws.Range("I73").Value = 4253
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4253
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3317
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 4571.2856
$ws.Range("I80").Value = 2999.6667
$ws.Range("J80").Value = 5750
$ws.Range("K80").Value = 2999.6667
$ws.Range("L80").Value = 5750
$ws.Range("M80").Value = -2001.6667
$ws.Range("N80").Value = -7746

$ws.Range("H83").Value = 4571.2856
$ws.Range("I83").Value = 2999.6667
$ws.Range("J83").Value = 5750
$ws.Range("K83").Value = 14998.3335
$ws.Range("L83").Value = 28750
$ws.Range("M83").Value = -10006.3335
$ws.Range("N83").Value = -38734

$ws.Range("H132").Value = 1499.6666
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7198.8
$ws.Range("I7").Value = 7427
$ws.Range("K7").Value = 7427
$ws.Range("M7").Value = -7315

$ws.Range("H40").Value = 5187.625
$ws.Range("I40").Value = 5000.1665
$ws.Range("K40").Value = 5000.1665
$ws.Range("M40").Value = -4864.1665

$ws.Range("H55").Value = 397.8
$ws.Range("I55").Value = 373.6
$ws.Range("J55").Value = 446.2
$ws.Range("K55").Value = 373.6
$ws.Range("L55").Value = 446.2
$ws.Range("M55").Value = -200.6
$ws.Range("N55").Value = -792.2

$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H122").Value = 5032
$ws.Range("I122").Value = 3451.8
$ws.Range("K122").Value = 10355.4
$ws.Range("M122").Value = -7905.400000000001

$ws.Range("H126").Value = 7198.8
$ws.Range("I126").Value = 7427
$ws.Range("K126").Value = 22281
$ws.Range("M126").Value = -19811

$ws.Range("H136").Value = 6308.3076
$ws.Range("I136").Value = 5855.875
$ws.Range("J136").Value = 7032.2
$ws.Range("K136").Value = 17567.625
$ws.Range("L136").Value = 21096.6
$ws.Range("M136").Value = -15017.625
$ws.Range("N136").Value = -26196.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 213395
$ws.Range("I49").Value = 229494.5
$ws.Range("K49").Value = 229494.5
$ws.Range("M49").Value = -229264.5

$ws.Range("H54").Value = 27898.2
$ws.Range("I54").Value = 27510
$ws.Range("J54").Value = 27995.25
$ws.Range("K54").Value = 27510
$ws.Range("L54").Value = 27995.25
$ws.Range("M54").Value = -26990
$ws.Range("N54").Value = -29035.25

$ws.Range("H81").Value = 4142.375
$ws.Range("I81").Value = 4142.375
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8284.75
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -7223.75
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 4142.375
$ws.Range("I84").Value = 4142.375
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 41423.75
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -36119.75
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 2777.111
$ws.Range("I113").Value = 724.1667
$ws.Range("J113").Value = 6883
$ws.Range("K113").Value = 2172.5001
$ws.Range("L113").Value = 20649
$ws.Range("M113").Value = -2.500100000000202
$ws.Range("N113").Value = -24989

$ws.Range("H122").Value = 4108.8
$ws.Range("I122").Value = 4116.643
$ws.Range("K122").Value = 12349.929
$ws.Range("M122").Value = -9899.929
